$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    25 = @{ A = 112183146; B = 96735; Q = 572346; R = 6714917 }
    26 = @{ A = 112183148; B = 96735; Q = 572357; R = 6714903 }
    27 = @{ A = 112183141; B = 96735; Q = 572361; R = 6714980 }
    28 = @{ A = 112183143; B = 96735; Q = 572359; R = 6714905 }
    29 = @{ A = 112183140; B = 96735; Q = 572350; R = 6714962 }
    30 = @{ A = 112183151; B = 96735; Q = 572361; R = 6714978 }
    31 = @{ A = 112183145; B = 96735; Q = 572351; R = 6714907 }
    32 = @{ A = 112183134; B = 96735; Q = 572354; R = 6714968 }
    33 = @{ A = 112183137; B = 96735; Q = 572354; R = 6714961 }
    34 = @{ A = 112183147; B = 96735; Q = 572351; R = 6714915 }
    35 = @{ A = 112183149; B = 96735; Q = 572345; R = 6714965 }
    36 = @{ A = 112183150; B = 96735; Q = 572358; R = 6714972 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
